$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 9).Value = 0.8781048434890718   # I2
$ws.Cells.Item(2, 10).Value = 0.8781048434890719  # J2
$ws.Cells.Item(2, 13).Value = 0.478362             # M2
$ws.Cells.Item(2, 14).Value = 1.435086             # N2
$ws.Cells.Item(2, 15).Value = 0.6533545125880439   # O2
$ws.Cells.Item(2, 16).Value = 0.6533545125880439   # P2
$ws.Cells.Item(2, 17).Value = 1.68806614916        # Q2
$ws.Cells.Item(2, 18).Value = 15.19259534244       # R2
$ws.Cells.Item(2, 19).Value = 0.573713762019003    # S2
$ws.Cells.Item(2, 20).Value = 0.5737137620190031   # T2

# Row 3
$ws.Cells.Item(3, 9).Value = 0.8781048434890718    # I3
$ws.Cells.Item(3, 10).Value = 0.8781048434890719   # J3
$ws.Cells.Item(3, 15).Value = 0.2124690813384451   # O3
$ws.Cells.Item(3, 16).Value = 0.2124690813384451   # P3
$ws.Cells.Item(3, 19).Value = 0.1865701294149622   # S3
$ws.Cells.Item(3, 20).Value = 0.1865701294149622   # T3

# Row 4
$ws.Cells.Item(4, 9).Value = 0.8781048434890718    # I4
$ws.Cells.Item(4, 10).Value = 0.8781048434890719   # J4
$ws.Cells.Item(4, 13).Value = 0.09823900000000001  # M4
$ws.Cells.Item(4, 14).Value = 0.294717              # N4
$ws.Cells.Item(4, 15).Value = 0.134176406073511    # O4
$ws.Cells.Item(4, 16).Value = 0.1341764060735109   # P4
$ws.Cells.Item(4, 17).Value = 0.3466703676866666   # Q4
$ws.Cells.Item(4, 18).Value = 3.12003330918        # R4
$ws.Cells.Item(4, 19).Value = 0.1178209520551065   # S4
$ws.Cells.Item(4, 20).Value = 0.1178209520551065   # T4

# Row 5
$ws.Cells.Item(5, 5).Value = 3                      # E5
$ws.Cells.Item(5, 6).Value = 1                      # F5
$ws.Cells.Item(5, 7).Value = 0.489861               # G5
$ws.Cells.Item(5, 8).Value = 1.469583               # H5
$ws.Cells.Item(5, 9).Value = 0.1218951565109281    # I5
$ws.Cells.Item(5, 10).Value = 0.1218951565109281   # J5
$ws.Cells.Item(5, 13).Value = 0.478362              # M5
$ws.Cells.Item(5, 14).Value = 1.435086              # N5
$ws.Cells.Item(5, 15).Value = 0.6533545125880439   # O5
$ws.Cells.Item(5, 16).Value = 0.6533545125880439   # P5
$ws.Cells.Item(5, 17).Value = 0.234330887682        # Q5
$ws.Cells.Item(5, 18).Value = 2.108977989138        # R5
$ws.Cells.Item(5, 19).Value = 0.07964075056904076  # S5
$ws.Cells.Item(5, 20).Value = 0.07964075056904076  # T5

# Row 6
$ws.Cells.Item(6, 5).Value = 3                      # E6
$ws.Cells.Item(6, 6).Value = 1                      # F6
$ws.Cells.Item(6, 7).Value = 0.489861               # G6
$ws.Cells.Item(6, 8).Value = 1.469583               # H6
$ws.Cells.Item(6, 9).Value = 0.1218951565109281    # I6
$ws.Cells.Item(6, 10).Value = 0.1218951565109281   # J6
$ws.Cells.Item(6, 15).Value = 0.2124690813384451   # O6
$ws.Cells.Item(6, 16).Value = 0.2124690813384451   # P6
$ws.Cells.Item(6, 17).Value = 0.07620375688200001  # Q6
$ws.Cells.Item(6, 18).Value = 0.685833811938        # R6
$ws.Cells.Item(6, 19).Value = 0.02589895192348288  # S6
$ws.Cells.Item(6, 20).Value = 0.02589895192348288  # T6

# Row 7
$ws.Cells.Item(7, 5).Value = 3                      # E7
$ws.Cells.Item(7, 6).Value = 1                      # F7
$ws.Cells.Item(7, 7).Value = 0.489861               # G7
$ws.Cells.Item(7, 8).Value = 1.469583               # H7
$ws.Cells.Item(7, 9).Value = 0.1218951565109281    # I7
$ws.Cells.Item(7, 10).Value = 0.1218951565109281   # J7
$ws.Cells.Item(7, 13).Value = 0.09823900000000001  # M7
$ws.Cells.Item(7, 14).Value = 0.294717              # N7
$ws.Cells.Item(7, 15).Value = 0.134176406073511    # O7
$ws.Cells.Item(7, 16).Value = 0.1341764060735109   # P7
$ws.Cells.Item(7, 17).Value = 0.04812345477900001  # Q7
$ws.Cells.Item(7, 18).Value = 0.433111093011        # R7
$ws.Cells.Item(7, 19).Value = 0.01635545401840446  # S7
$ws.Cells.Item(7, 20).Value = 0.01635545401840446  # T7
